$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing datetime number format from A2 before we touch it,
# so the new A3 cell can reuse that same format.
$dateFormat = $ws.Range("A2").NumberFormat

# A1: 42 -> 52879 (stays a plain general number, no style change)
$ws.Range("A1").Value = 52879

# A2: was a date-formatted value (44987.66474087031) -> becomes plain 1.
# Fully clear its formatting (drops the date-time style back to default)
# before writing the new plain number.
$ws.Range("A2").ClearFormats()
$ws.Range("A2").Value = 1

# New row 3: A3 gets the date-time value/format that used to live in A2
$ws.Range("A3").NumberFormat = $dateFormat
$ws.Range("A3").Value = 44987.76033516775
